# Applies the "feat: add 2022-Q3 data" edit:
#  1. Inserts a new "2022-Q3" worksheet right after "总计", populated with a
#     copy of the "2022-Q2" sheet's layout/styling but with the new quarter's
#     fund data.
#  2. Updates the "总计" (summary) sheet: inserts a new row at the top of the
#     data (shifting the existing quarters down by one row) with the new
#     2022-Q3 summary figures.
#  3. Fixes a mislabeled header ("基金金额" -> "基金规模") on the "2021-Q3"
#     sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by duplicating "2022-Q2" (keeps styles),
#    placed immediately before "2022-Q2" (i.e. right after "总计").
# ---------------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Copy($q2Sheet, [System.Reflection.Missing]::Value)
$q3Sheet = $wb.Worksheets.Item("2022-Q2 (2)")
$q3Sheet.Name = "2022-Q3"

# Row 2: 002423 / 华宝标普美国品质消费股票（LOF）美元
$q3Sheet.Range("C2").Value = "华宝标普美国品质消费股票（LOF）美元"
$q3Sheet.Range("D2:G2").NumberFormat = "@"
$q3Sheet.Range("D2").Value = "3.59"
$q3Sheet.Range("E2").Value = "94.37"
$q3Sheet.Range("F2").Value = "2.09"
$q3Sheet.Range("G2").Value = "0.0750"
$q3Sheet.Range("H2").Value = 10

# Row 3: 162415 / 华宝标普美国品质消费股票（LOF）人民币A
$q3Sheet.Range("C3").Value = "华宝标普美国品质消费股票（LOF）人民币A"
$q3Sheet.Range("D3:G3").NumberFormat = "@"
$q3Sheet.Range("D3").Value = "2.86"
$q3Sheet.Range("E3").Value = "94.37"
$q3Sheet.Range("F3").Value = "2.09"
$q3Sheet.Range("G3").Value = "0.0598"
$q3Sheet.Range("H3").Value = 10

# Row 4: 009975 / 华宝标普美国品质消费股票（LOF）人民币C
$q3Sheet.Range("C4").Value = "华宝标普美国品质消费股票（LOF）人民币C"
$q3Sheet.Range("D4:G4").NumberFormat = "@"
$q3Sheet.Range("D4").Value = "0.73"
$q3Sheet.Range("E4").Value = "94.37"
$q3Sheet.Range("F4").Value = "2.09"
$q3Sheet.Range("G4").Value = "0.0153"
$q3Sheet.Range("H4").Value = 10

# ---------------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new data row right under the header
#    and fill it with the 2022-Q3 summary; the other rows shift down
#    automatically (Excel keeps their own content untouched).
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows("2:2").Insert()

# Match formatting of the row below (which holds the original 2022-Q2 row).
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2:D2").PasteSpecial(-4122)

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 3
$totalSheet.Cells.Item(2, 4).Value = 0.15

# ---------------------------------------------------------------------------
# 3. Fix the "基金金额" -> "基金规模" header typo on the "2021-Q3" sheet.
# ---------------------------------------------------------------------------
$q3_2021Sheet = $wb.Worksheets.Item("2021-Q3")
$q3_2021Sheet.Range("D1").Value = "基金规模"
